# Check-URL-Redirects/Inputs/Example.xlsx edit
# - rename sheet Sheet1 -> Redirects
# - update the saved absPath / window position
# - upgrade a handful of "http://" targets (that are bare-domain redirects) to "https://"
# - hyperlink those same cells to their https target
# - drop the unused conditional-format-ish styles, keep only Normal + bold header

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Redirects"

$win = $excel.ActiveWindow
$win.Left = 1200
$win.Top = 1300

# URLs that only had a bare host (no path) get bumped from http to https and
# turned into real hyperlinks; URLs that already pointed at a full path were
# left alone per the commit message.
$updates = @(
    @{ Cell = "B2";  Url = "https://www.masterlock.com" },
    @{ Cell = "B3";  Url = "https://www.masterlock.comz" },
    @{ Cell = "B4";  Url = "https://www.masterlock.com/cms/customersupport/current-pricing" },
    @{ Cell = "B7";  Url = "https://www.masterlock.com/service-and-support/faqs/lost-combinations" },
    @{ Cell = "B13"; Url = "https://www.masterlock.com/cms/customersupport/current-pricing" },
    @{ Cell = "B16"; Url = "https://www.masterlock.com/personal-use/product/175DLH" },
    @{ Cell = "B17"; Url = "https://www.masterlock.com/business-use/product-search/safety-solutions" },
    @{ Cell = "B18"; Url = "https://www.masterlock.com/business-use/commercial-security" },
    @{ Cell = "B19"; Url = "https://www.masterlock.com/legal-statement" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.Value = $u.Url
    $ws.Hyperlinks.Add($rng, $u.Url)
}

# Clear stray "Good"/"Bad" style formatting picked up by earlier edits of this
# sheet - only the bold header row keeps explicit formatting.
$ws.Range("A2:B19").ClearFormats()

$ws.Range("A21").Select()

Write-Output "done"
